# Casos de prueba actualizados
# Se definieron los valores a ingresar en los campos
#
# Replace the generic <Placeholder> tokens in the 'Precondiciones' and 'Pasos'
# sheets with the concrete values that will actually be entered when the test
# case is executed (CPA_Playa1 / CPA_Playa2 / fixed phone, mail, schedule, etc.).

$wb = $excel.ActiveWorkbook
$wsDatos = $wb.Worksheets.Item("DatosGenerales")
$wsPre   = $wb.Worksheets.Item("Precondiciones")
$wsPasos = $wb.Worksheets.Item("Pasos")
$wsCtrl  = $wb.Worksheets.Item("Control de cambios")

$wsPre.Range("B3").Value = '"CPA_Playa1" es el nombre de la playa <Playa1>'
$wsPre.Range("B4").Value = ' "999999" es un telefono valido de playa de estacionamiento'
$wsPre.Range("B5").Value = ' "CPA_mail@CPA_mail" es un mail válido de la playa de estacionamiento'
$wsPre.Range("B6").Value = '<CPA_TipoPlaya1> es un tipo de playa de estacionamiento valido y existe en la base de datos'
$wsPre.Range("B7").Value = ' <CPA_TipoVehiculo1> es un tipo de vehiculo que existe en la base de datos y <capacidad1> es un valor válido para una capacidad'
$wsPre.Range("B8").Value = ' <CPA_TipoVehiculo2> es un tipo de vehiculo que existe en la base de datos y <Capacidad2> es un valor válido para una capacidad'
$wsPre.Range("B9").Value = '<Domicilio1>, conla provincia <Córdoba> existe en la base de datos,  el departamento <Capital> existe en la base de datos y pertenece a  <Córdoba>, la ciudad <Córdoba> existe en la base de datos y pertenece al <Capital>, <Calle1> con nombre "Colon" es una calle valida y "9" es un valor valido para numero; es domicilio valido para una playa de estacionamiento'
$wsPre.Range("B10").Value = '<Horario1> con <CPA_DiasDeAtencion1> como dias, "00:00" como horario desde y "23:59" hasta es un horario de <Playa1>'
$wsPre.Range("B11").Value = ' <Precio1> con <CPA_TipoVehiculo1> como tipo de vehiculo, <CPA_TipoHorario1> como tipo de horario, <CPA_DiasDeAtencion1> como dias y "9" como precio es un precio de <Playa1>'
$wsPre.Range("B12").Value = '"CPA_Playa2" es un nombre válido para una playa de estcionamiento'
$wsPasos.Range("B3").Value = 'Ingreso "CPA_Playa1" en el campo nombre de playa'
$wsPasos.Range("B6").Value = 'ingreso "CPA_Playa2" en el campo nombre de playa'
$wsPasos.Range("C7").Value = 'Se muestra el mensaje "Está seguro que desea guardar los cambios en la playa de estacionamiento CPA_Playa1?"'

# Restore the per-sheet selection and make 'Precondiciones' the active tab,
# matching where the author was working when the values were filled in.
$wsDatos.Range("B4").Select()
$wsPasos.Range("C6").Select()
$wsCtrl.Range("D3").Select()
$wsPre.Range("B6").Select()
$wsPre.Activate()

